$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting for the two newly appended rows (48, 49) from row 47
$ws.Range("A47:F47").Copy()
$ws.Range("A48:F49").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 43271
$ws.Cells.Item(2, 3).Value = 1054
$ws.Cells.Item(2, 4).Value = 13.68181818181818
$ws.Cells.Item(2, 5).Value = "Anta-Espinho"
$ws.Cells.Item(2, 6).Value = 20.64603079546817

# Row 3
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 43271
$ws.Cells.Item(3, 3).Value = 3093
$ws.Cells.Item(3, 4).Value = 36.49166666666667
$ws.Cells.Item(3, 5).Value = "Arcos"
$ws.Cells.Item(3, 6).Value = 32.93764130559803

# Row 4
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 43271
$ws.Cells.Item(4, 3).Value = 2017
$ws.Cells.Item(4, 4).Value = 27.125
$ws.Cells.Item(4, 5).Value = "Aveiro"
$ws.Cells.Item(4, 6).Value = 18.85869837746607

# Row 5
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 43271
$ws.Cells.Item(5, 3).Value = 3075
$ws.Cells.Item(5, 4).Value = 48.41666666666666
$ws.Cells.Item(5, 5).Value = "Avenida da Liberdade"
$ws.Cells.Item(5, 6).Value = 32.89084930364226

# Row 6
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 43271
$ws.Cells.Item(6, 3).Value = 1053
$ws.Cells.Item(6, 4).Value = 27.2
$ws.Cells.Item(6, 5).Value = "Avintes"
$ws.Cells.Item(6, 6).Value = 21.89124717822572

# Row 7
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 43271
$ws.Cells.Item(7, 3).Value = 1052
$ws.Cells.Item(7, 4).Value = 24.45454545454545
$ws.Cells.Item(7, 5).Value = "Burgães-Santo Tirso"
$ws.Cells.Item(7, 6).Value = 18.73569299496276

# Row 8
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 43271
$ws.Cells.Item(8, 3).Value = 3104
$ws.Cells.Item(8, 4).Value = 41.12380952380953
$ws.Cells.Item(8, 5).Value = "Cascais - Escola da Cidadela"
$ws.Cells.Item(8, 6).Value = 28.20358804427715

# Row 9
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 43271
$ws.Cells.Item(9, 3).Value = 5012
$ws.Cells.Item(9, 4).Value = 14.45833333333333
$ws.Cells.Item(9, 5).Value = "Cerro"
$ws.Cells.Item(9, 6).Value = 16.79738220759935

# Row 10
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 43271
$ws.Cells.Item(10, 3).Value = 3096
$ws.Cells.Item(10, 4).Value = 29.04583333333333
$ws.Cells.Item(10, 5).Value = "Chamusca"
$ws.Cells.Item(10, 6).Value = 23.06041636353474

# Row 11
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 43271
$ws.Cells.Item(11, 3).Value = 2006
$ws.Cells.Item(11, 4).Value = 37.58333333333334
$ws.Cells.Item(11, 5).Value = "Coimbra/ Avenida Fernão Magalhães"
$ws.Cells.Item(11, 6).Value = 24.21703155113064

# Row 12
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 43271
$ws.Cells.Item(12, 3).Value = 1046
$ws.Cells.Item(12, 4).Value = 18.83333333333333
$ws.Cells.Item(12, 5).Value = "Cónego Dr. Manuel Faria-Azurém"
$ws.Cells.Item(12, 6).Value = 15.45705537044591

# Row 13
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 43271
$ws.Cells.Item(13, 3).Value = 5011
$ws.Cells.Item(13, 4).Value = 38.29166666666666
$ws.Cells.Item(13, 5).Value = "David Neto"
$ws.Cells.Item(13, 6).Value = 29.74194464589739

# Row 14
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 43271
$ws.Cells.Item(14, 3).Value = 1048
$ws.Cells.Item(14, 4).Value = 10.75
$ws.Cells.Item(14, 5).Value = "Douro Norte"
$ws.Cells.Item(14, 6).Value = 9.911897496874564

# Row 15
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 43271
$ws.Cells.Item(15, 3).Value = 3072
$ws.Cells.Item(15, 4).Value = 40.19583333333333
$ws.Cells.Item(15, 5).Value = "Entrecampos"
$ws.Cells.Item(15, 6).Value = 31.7736527854409

# Row 16
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 43271
$ws.Cells.Item(16, 3).Value = 1023
$ws.Cells.Item(16, 4).Value = 27.20833333333333
$ws.Cells.Item(16, 5).Value = "Ermesinde-Valongo"
$ws.Cells.Item(16, 6).Value = 22.68205445900004

# Row 17
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 43271
$ws.Cells.Item(17, 3).Value = 2019
$ws.Cells.Item(17, 4).Value = 31.45833333333333
$ws.Cells.Item(17, 5).Value = "Ervedeira"
$ws.Cells.Item(17, 6).Value = 25.17787668475592

# Row 18
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 43271
$ws.Cells.Item(18, 3).Value = 3095
$ws.Cells.Item(18, 4).Value = 35.2875
$ws.Cells.Item(18, 5).Value = "Escavadeira"
$ws.Cells.Item(18, 6).Value = 33.82417065389357

# Row 19
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = 43271
$ws.Cells.Item(19, 3).Value = 2004
$ws.Cells.Item(19, 4).Value = 17.89473684210526
$ws.Cells.Item(19, 5).Value = "Estarreja"
$ws.Cells.Item(19, 6).Value = 18.94870887804427

# Row 20
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = 43271
$ws.Cells.Item(20, 3).Value = 3099
$ws.Cells.Item(20, 4).Value = 38.3375
$ws.Cells.Item(20, 5).Value = "Fernando Pó"
$ws.Cells.Item(20, 6).Value = 28.03150538351018

# Row 21
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = 43271
$ws.Cells.Item(21, 3).Value = 2021
$ws.Cells.Item(21, 4).Value = 14.47368421052632
$ws.Cells.Item(21, 5).Value = "Fornelo do Monte"
$ws.Cells.Item(21, 6).Value = 14.67341021479265

# Row 22
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = 43271
$ws.Cells.Item(22, 3).Value = 1028
$ws.Cells.Item(22, 4).Value = 17.22727272727273
$ws.Cells.Item(22, 5).Value = "Francisco Sá Carneiro-Campanha"
$ws.Cells.Item(22, 6).Value = 23.03242851249266

# Row 23
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = 43271
$ws.Cells.Item(23, 3).Value = 1042
$ws.Cells.Item(23, 4).Value = 10.08333333333333
$ws.Cells.Item(23, 5).Value = "Frossos-Braga"
$ws.Cells.Item(23, 6).Value = 15.10060868609654

# Row 24
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = 43271
$ws.Cells.Item(24, 3).Value = 2020
$ws.Cells.Item(24, 4).Value = 18.41666666666667
$ws.Cells.Item(24, 5).Value = "Fundão"
$ws.Cells.Item(24, 6).Value = 14.6387783591567

# Row 25
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = 43271
$ws.Cells.Item(25, 3).Value = 2018
$ws.Cells.Item(25, 4).Value = 32.83333333333334
$ws.Cells.Item(25, 5).Value = "Ílhavo"
$ws.Cells.Item(25, 6).Value = 18.55732009616839

# Row 26
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = 43271
$ws.Cells.Item(26, 3).Value = 2016
$ws.Cells.Item(26, 4).Value = 28.45833333333333
$ws.Cells.Item(26, 5).Value = "Instituto Geofísico de Coimbra"
$ws.Cells.Item(26, 6).Value = 24.40424949410358

# Row 27
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = 43271
$ws.Cells.Item(27, 3).Value = 1030
$ws.Cells.Item(27, 4).Value = 27.28571428571428
$ws.Cells.Item(27, 5).Value = "João Gomes Laranjo-S.Hora"
$ws.Cells.Item(27, 6).Value = 23.26173016194345

# Row 28
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = 43271
$ws.Cells.Item(28, 3).Value = 5007
$ws.Cells.Item(28, 4).Value = 30.44166666666667
$ws.Cells.Item(28, 5).Value = "Joaquim Magalhães"
$ws.Cells.Item(28, 6).Value = 27.02458913103547

# Row 29
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = 43271
$ws.Cells.Item(29, 3).Value = 3083
$ws.Cells.Item(29, 4).Value = 38.62916666666667
$ws.Cells.Item(29, 5).Value = "Laranjeiro"
$ws.Cells.Item(29, 6).Value = 32.97233577505904

# Row 30
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = 43271
$ws.Cells.Item(30, 3).Value = 3055
$ws.Cells.Item(30, 4).Value = 39.97916666666666
$ws.Cells.Item(30, 5).Value = "Lavradio"
$ws.Cells.Item(30, 6).Value = 33.6246065931518

# Row 31
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = 43271
$ws.Cells.Item(31, 3).Value = 3085
$ws.Cells.Item(31, 4).Value = 25.6375
$ws.Cells.Item(31, 5).Value = "Loures-Centro"
$ws.Cells.Item(31, 6).Value = 27.97678762543504

# Row 32
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = 43271
$ws.Cells.Item(32, 3).Value = 3102
$ws.Cells.Item(32, 4).Value = 25.84583333333333
$ws.Cells.Item(32, 5).Value = "Lourinhã"
$ws.Cells.Item(32, 6).Value = 25.45935860930451

# Row 33
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).Value = 43271
$ws.Cells.Item(33, 3).Value = 5008
$ws.Cells.Item(33, 4).Value = 43.41666666666666
$ws.Cells.Item(33, 5).Value = "Malpique"
$ws.Cells.Item(33, 6).Value = 29.89094651229147

# Row 34
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = 43271
$ws.Cells.Item(34, 3).Value = 1025
$ws.Cells.Item(34, 4).Value = 20.25
$ws.Cells.Item(34, 5).Value = "Meco-Perafita"
$ws.Cells.Item(34, 6).Value = 22.66369613285387

# Row 35
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = 43271
$ws.Cells.Item(35, 3).Value = 3089
$ws.Cells.Item(35, 4).Value = 27.72083333333333
$ws.Cells.Item(35, 5).Value = "Mem Martins"
$ws.Cells.Item(35, 6).Value = 26.55954726445246

# Row 36
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = 43271
$ws.Cells.Item(36, 3).Value = 1051
$ws.Cells.Item(36, 4).Value = 16.58333333333333
$ws.Cells.Item(36, 5).Value = "Mindelo-Vila do Conde"
$ws.Cells.Item(36, 6).Value = 21.24585750058139

# Row 37
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = 43271
$ws.Cells.Item(37, 3).Value = 2022
$ws.Cells.Item(37, 4).Value = 28.625
$ws.Cells.Item(37, 5).Value = "Montemor-o-Velho"
$ws.Cells.Item(37, 6).Value = 25.25524961832934

# Row 38
$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(38, 2).Value = 43271
$ws.Cells.Item(38, 3).Value = 3097
$ws.Cells.Item(38, 4).Value = 43.25416666666666
$ws.Cells.Item(38, 5).Value = "Odivelas-Ramada"
$ws.Cells.Item(38, 6).Value = 28.63738323879013

# Row 39
$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(39, 2).Value = 43271
$ws.Cells.Item(39, 3).Value = 3071
$ws.Cells.Item(39, 4).Value = 36.34166666666667
$ws.Cells.Item(39, 5).Value = "Olivais"
$ws.Cells.Item(39, 6).Value = 31.03275876896798

# Row 40
$ws.Cells.Item(40, 1).Value = 38
$ws.Cells.Item(40, 2).Value = 43271
$ws.Cells.Item(40, 3).Value = 1044
$ws.Cells.Item(40, 4).Value = 27.79166666666667
$ws.Cells.Item(40, 5).Value = "Paços de Ferreira"
$ws.Cells.Item(40, 6).Value = 18.62189344545493

# Row 41
$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 2).Value = 43271
$ws.Cells.Item(41, 3).Value = 3063
$ws.Cells.Item(41, 4).Value = 38.65833333333333
$ws.Cells.Item(41, 5).Value = "Paio Pires"
$ws.Cells.Item(41, 6).Value = 33.57066858262246

# Row 42
$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 2).Value = 43271
$ws.Cells.Item(42, 3).Value = 1043
$ws.Cells.Item(42, 4).Value = 7
$ws.Cells.Item(42, 5).Value = "Pe Moreira Neves-Castelões de Cepeda"
$ws.Cells.Item(42, 6).Value = 18.08886559205197

# Row 43
$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 2).Value = 43271
$ws.Cells.Item(43, 3).Value = 3094
$ws.Cells.Item(43, 4).Value = 39.85833333333333
$ws.Cells.Item(43, 5).Value = "Quebedo"
$ws.Cells.Item(43, 6).Value = 33.01948570155671

# Row 44
$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 2).Value = 43271
$ws.Cells.Item(44, 3).Value = 3091
$ws.Cells.Item(44, 4).Value = 28.62916666666667
$ws.Cells.Item(44, 5).Value = "Quinta do Marquês"
$ws.Cells.Item(44, 6).Value = 29.65253984794291

# Row 45
$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(45, 2).Value = 43271
$ws.Cells.Item(45, 3).Value = 3084
$ws.Cells.Item(45, 4).Value = 31.62083333333333
$ws.Cells.Item(45, 5).Value = "Reboleira"
$ws.Cells.Item(45, 6).Value = 29.48644693160862

# Row 46
$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = 43271
$ws.Cells.Item(46, 3).Value = 1055
$ws.Cells.Item(46, 4).Value = 25
$ws.Cells.Item(46, 5).Value = "Seara-Matosinhos"
$ws.Cells.Item(46, 6).Value = 23.30881037668967

# Row 47
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).Value = 43271
$ws.Cells.Item(47, 3).Value = 4003
$ws.Cells.Item(47, 4).Value = 19.29166666666667
$ws.Cells.Item(47, 5).Value = "Sonega"
$ws.Cells.Item(47, 6).Value = 24.91582792010203

# Row 48
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = 43271
$ws.Cells.Item(48, 3).Value = 4006
$ws.Cells.Item(48, 4).Value = 25.33333333333333
$ws.Cells.Item(48, 5).Value = "Terena"
$ws.Cells.Item(48, 6).Value = 22.57544372407552

# Row 49
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).Value = 43271
$ws.Cells.Item(49, 3).Value = 1031
$ws.Cells.Item(49, 4).Value = 25.70588235294118
$ws.Cells.Item(49, 5).Value = "VNTelha-Maia"
$ws.Cells.Item(49, 6).Value = 22.30368022065842
